$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 303
$ws.Range("A303").NumberFormat = "@"
$ws.Cells.Item(303,1).Value = "MAYOR"
$ws.Range("A303").ClearFormats()

# Row 305
$ws.Range("A305:C305").NumberFormat = "@"
$ws.Cells.Item(305,1).Value = "Candidate"
$ws.Cells.Item(305,2).Value = "Votes"
$ws.Cells.Item(305,3).Value = "Percentage"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A305:C305").PasteSpecial(-4122) | Out-Null

# Row 306
$ws.Range("A306:C306").NumberFormat = "@"
$ws.Cells.Item(306,1).Value = "CAYETANO, ATE LANI (NP)"
$ws.Cells.Item(306,2).Value = "272,876"
$ws.Cells.Item(306,3).Value = "75.76 %"
$ws.Range("A306:C306").ClearFormats()

# Row 307
$ws.Range("A307:C307").NumberFormat = "@"
$ws.Cells.Item(307,1).Value = "CERAFICA, ARNEL (PPP)"
$ws.Cells.Item(307,2).Value = "87,266"
$ws.Cells.Item(307,3).Value = "24.23 %"
$ws.Range("A307:C307").ClearFormats()

# Row 309
$ws.Range("A309:D309").NumberFormat = "@"
$ws.Cells.Item(309,1).Value = "Over-votes"
$ws.Cells.Item(309,2).Value = "Under-votes"
$ws.Cells.Item(309,3).Value = "Valid votes"
$ws.Cells.Item(309,4).Value = "Votes obtained by all candidates"
$ws.Range("A289:D289").Copy() | Out-Null
$ws.Range("A309:D309").PasteSpecial(-4122) | Out-Null

# Row 310
$ws.Range("A310:D310").NumberFormat = "@"
$ws.Cells.Item(310,1).Value = "1183"
$ws.Cells.Item(310,2).Value = "9723"
$ws.Cells.Item(310,3).Value = "371575"
$ws.Cells.Item(310,4).Value = "360142"
$ws.Range("A310:D310").ClearFormats()

# Row 313
$ws.Range("A313").NumberFormat = "@"
$ws.Cells.Item(313,1).Value = "VICE-MAYOR"
$ws.Range("A313").ClearFormats()

# Row 315
$ws.Range("A315:C315").NumberFormat = "@"
$ws.Cells.Item(315,1).Value = "Candidate"
$ws.Cells.Item(315,2).Value = "Votes"
$ws.Cells.Item(315,3).Value = "Percentage"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A315:C315").PasteSpecial(-4122) | Out-Null

# Row 316
$ws.Range("A316:C316").NumberFormat = "@"
$ws.Cells.Item(316,1).Value = "ALIT, ARVIN (NP)"
$ws.Cells.Item(316,2).Value = "232,034"
$ws.Cells.Item(316,3).Value = "68.81 %"
$ws.Range("A316:C316").ClearFormats()

# Row 317
$ws.Range("A317:C317").NumberFormat = "@"
$ws.Cells.Item(317,1).Value = "CERAFICA, JANELLE (PPP)"
$ws.Cells.Item(317,2).Value = "105,157"
$ws.Cells.Item(317,3).Value = "31.18 %"
$ws.Range("A317:C317").ClearFormats()

# Row 319
$ws.Range("A319:D319").NumberFormat = "@"
$ws.Cells.Item(319,1).Value = "Over-votes"
$ws.Cells.Item(319,2).Value = "Under-votes"
$ws.Cells.Item(319,3).Value = "Valid votes"
$ws.Cells.Item(319,4).Value = "Votes obtained by all candidates"
$ws.Range("A289:D289").Copy() | Out-Null
$ws.Range("A319:D319").PasteSpecial(-4122) | Out-Null

# Row 320
$ws.Range("A320:D320").NumberFormat = "@"
$ws.Cells.Item(320,1).Value = "384"
$ws.Cells.Item(320,2).Value = "33473"
$ws.Cells.Item(320,3).Value = "371575"
$ws.Cells.Item(320,4).Value = "337191"
$ws.Range("A320:D320").ClearFormats()

# Row 323
$ws.Range("A323").NumberFormat = "@"
$ws.Cells.Item(323,1).Value = "MEMBER, SANGGUNIANG PANLUNGSOD"
$ws.Range("A323").ClearFormats()

# Row 325
$ws.Range("A325:C325").NumberFormat = "@"
$ws.Cells.Item(325,1).Value = "Candidate"
$ws.Cells.Item(325,2).Value = "Votes"
$ws.Cells.Item(325,3).Value = "Percentage"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A325:C325").PasteSpecial(-4122) | Out-Null

# Row 326
$ws.Range("A326:C326").NumberFormat = "@"
$ws.Cells.Item(326,1).Value = "AMOROSO, INOCENTES (IND)"
$ws.Cells.Item(326,2).Value = "10,681"
$ws.Cells.Item(326,3).Value = "0.97 %"
$ws.Range("A326:C326").ClearFormats()

# Row 327
$ws.Range("A327:C327").NumberFormat = "@"
$ws.Cells.Item(327,1).Value = "AQUINO, COMMISSIONER (NP)"
$ws.Cells.Item(327,2).Value = "85,159"
$ws.Cells.Item(327,3).Value = "7.77 %"
$ws.Range("A327:C327").ClearFormats()

# Row 328
$ws.Range("A328:C328").NumberFormat = "@"
$ws.Cells.Item(328,1).Value = "BAUTISTA, JONJON (PPP)"
$ws.Cells.Item(328,2).Value = "33,826"
$ws.Cells.Item(328,3).Value = "3.08 %"
$ws.Range("A328:C328").ClearFormats()

# Row 329
$ws.Range("A329:C329").NumberFormat = "@"
$ws.Cells.Item(329,1).Value = "BERNAL, RJ (PPP)"
$ws.Cells.Item(329,2).Value = "27,270"
$ws.Cells.Item(329,3).Value = "2.49 %"
$ws.Range("A329:C329").ClearFormats()

# Row 330
$ws.Range("A330:C330").NumberFormat = "@"
$ws.Cells.Item(330,1).Value = "CRUZ, JANNAH (PPP)"
$ws.Cells.Item(330,2).Value = "38,925"
$ws.Cells.Item(330,3).Value = "3.55 %"
$ws.Range("A330:C330").ClearFormats()

# Row 331
$ws.Range("A331:C331").NumberFormat = "@"
$ws.Cells.Item(331,1).Value = "DE MESA, GIGI VALENZUELA (NP)"
$ws.Cells.Item(331,2).Value = "105,617"
$ws.Cells.Item(331,3).Value = "9.64 %"
$ws.Range("A331:C331").ClearFormats()

# Row 332
$ws.Range("A332:C332").NumberFormat = "@"
$ws.Cells.Item(332,1).Value = "DELOS SANTOS, WARREN (PPP)"
$ws.Cells.Item(332,2).Value = "35,586"
$ws.Cells.Item(332,3).Value = "3.25 %"
$ws.Range("A332:C332").ClearFormats()

# Row 333
$ws.Range("A333:C333").NumberFormat = "@"
$ws.Cells.Item(333,1).Value = "DIONISIO, WARREN (PPP)"
$ws.Cells.Item(333,2).Value = "41,554"
$ws.Cells.Item(333,3).Value = "3.79 %"
$ws.Range("A333:C333").ClearFormats()

# Row 334
$ws.Range("A334:C334").NumberFormat = "@"
$ws.Cells.Item(334,1).Value = "FRANCO, RONET (PPP)"
$ws.Cells.Item(334,2).Value = "35,560"
$ws.Cells.Item(334,3).Value = "3.24 %"
$ws.Range("A334:C334").ClearFormats()

# Row 335
$ws.Range("A335:C335").NumberFormat = "@"
$ws.Cells.Item(335,1).Value = "ICAY, ANGGUS (PPP)"
$ws.Cells.Item(335,2).Value = "50,475"
$ws.Cells.Item(335,3).Value = "4.61 %"
$ws.Range("A335:C335").ClearFormats()

# Row 336
$ws.Range("A336:C336").NumberFormat = "@"
$ws.Cells.Item(336,1).Value = "LABAMPA, JIMMY (NP)"
$ws.Cells.Item(336,2).Value = "105,771"
$ws.Cells.Item(336,3).Value = "9.66 %"
$ws.Range("A336:C336").ClearFormats()

# Row 337
$ws.Range("A337:C337").NumberFormat = "@"
$ws.Cells.Item(337,1).Value = "LONTOC, PAUL (IND)"
$ws.Cells.Item(337,2).Value = "24,923"
$ws.Cells.Item(337,3).Value = "2.27 %"
$ws.Range("A337:C337").ClearFormats()

# Row 338
$ws.Range("A338:C338").NumberFormat = "@"
$ws.Cells.Item(338,1).Value = "MADRID, ELVIRA (PPP)"
$ws.Cells.Item(338,2).Value = "24,620"
$ws.Cells.Item(338,3).Value = "2.24 %"
$ws.Range("A338:C338").ClearFormats()

# Row 339
$ws.Range("A339:C339").NumberFormat = "@"
$ws.Cells.Item(339,1).Value = "MAÑOSCA, TOTONG (NP)"
$ws.Cells.Item(339,2).Value = "83,719"
$ws.Cells.Item(339,3).Value = "7.64 %"
$ws.Range("A339:C339").ClearFormats()

# Row 340
$ws.Range("A340:C340").NumberFormat = "@"
$ws.Cells.Item(340,1).Value = "MARCELINO, TIKBOY (NP)"
$ws.Cells.Item(340,2).Value = "96,695"
$ws.Cells.Item(340,3).Value = "8.83 %"
$ws.Range("A340:C340").ClearFormats()

# Row 341
$ws.Range("A341:C341").NumberFormat = "@"
$ws.Cells.Item(341,1).Value = "OGALINOLA, CARLITO (NP)"
$ws.Cells.Item(341,2).Value = "93,044"
$ws.Cells.Item(341,3).Value = "8.49 %"
$ws.Range("A341:C341").ClearFormats()

# Row 342
$ws.Range("A342:C342").NumberFormat = "@"
$ws.Cells.Item(342,1).Value = "OSORIO, MAR NORBERT (IND)"
$ws.Cells.Item(342,2).Value = "12,696"
$ws.Cells.Item(342,3).Value = "1.15 %"
$ws.Range("A342:C342").ClearFormats()

# Row 343
$ws.Range("A343:C343").NumberFormat = "@"
$ws.Cells.Item(343,1).Value = "PANGA-CRUZ, ATTYJOY (NP)"
$ws.Cells.Item(343,2).Value = "94,514"
$ws.Cells.Item(343,3).Value = "8.63 %"
$ws.Range("A343:C343").ClearFormats()

# Row 344
$ws.Range("A344:C344").NumberFormat = "@"
$ws.Cells.Item(344,1).Value = "SAN PEDRO, GAMIE (NP)"
$ws.Cells.Item(344,2).Value = "94,111"
$ws.Cells.Item(344,3).Value = "8.59 %"
$ws.Range("A344:C344").ClearFormats()

# Row 346
$ws.Range("A346:D346").NumberFormat = "@"
$ws.Cells.Item(346,1).Value = "Over-votes"
$ws.Cells.Item(346,2).Value = "Under-votes"
$ws.Cells.Item(346,3).Value = "Valid votes"
$ws.Cells.Item(346,4).Value = "Votes obtained by all candidates"
$ws.Range("A289:D289").Copy() | Out-Null
$ws.Range("A346:D346").PasteSpecial(-4122) | Out-Null

# Row 347
$ws.Range("A347:D347").NumberFormat = "@"
$ws.Cells.Item(347,1).Value = "2266"
$ws.Cells.Item(347,2).Value = "245486"
$ws.Cells.Item(347,3).Value = "1360352"
$ws.Cells.Item(347,4).Value = "1094746"
$ws.Range("A347:D347").ClearFormats()

$excel.CutCopyMode = $false